$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updated s_val data (save games filtered)
$ws.Range("B2").Value = 3.182878228561681
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 0.7127328510149897
$ws.Range("E2").Value = 0.4998867070740569
$ws.Range("G2").Value = 6.048734245549538

# Row 3 updated s_val data (save games filtered)
$ws.Range("B3").Value = 0.3464964993005633
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 0.1529057820181812
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("G3").Value = 2.652525447291612
